$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 10100.61
$ws.Range("B6").Value = 10069.39
$ws.Range("C6").Value = 107.89
$ws.Range("D6").Value = 108.22
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = 0.31
$ws.Range("G6").Value = 42613.766597222224
$ws.Range("H6").Value = $true

# Row 7
$ws.Range("A7").Value = 10104.65
$ws.Range("B7").Value = 10100.61
$ws.Range("C7").Value = 107.17
$ws.Range("D7").Value = 107.21
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.04
$ws.Range("G7").Value = 42614.673842592594
$ws.Range("H7").Value = $true

# Row 8
$ws.Range("A8").Value = 10102.629999999999
$ws.Range("B8").Value = 10104.65
$ws.Range("C8").Value = 107.04
$ws.Range("D8").Value = 107.02
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -0.02
$ws.Range("G8").Value = 42615.752916666665
$ws.Range("H8").Value = $false
